# Scheduled-runner market data refresh: updates cached price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job sheets.
# No formulas are involved - every touched cell holds a literal number.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 5172.2
$ws.Range("J70").Value = 6249.5713
$ws.Range("L70").Value = 18748.7139
$ws.Range("N70").Value = -19288.7139

$ws.Range("H73").Value = 5172.2
$ws.Range("J73").Value = 6249.5713
$ws.Range("L73").Value = 18748.7139
$ws.Range("N73").Value = -20620.7139

$ws.Range("H80").Value = 202.88889
$ws.Range("I80").Value = 89.333336
$ws.Range("J80").Value = 430
$ws.Range("K80").Value = 268.000008
$ws.Range("L80").Value = 1290
$ws.Range("M80").Value = 729.999992
$ws.Range("N80").Value = -3286

$ws.Range("H83").Value = 202.88889
$ws.Range("I83").Value = 89.333336
$ws.Range("J83").Value = 430
$ws.Range("K83").Value = 804.0000240000001
$ws.Range("L83").Value = 3870
$ws.Range("M83").Value = 4187.999976
$ws.Range("N83").Value = -13854

$ws.Range("H86").Value = 1866.6666
$ws.Range("J86").Value = 2000
$ws.Range("L86").Value = 2000
$ws.Range("N86").Value = -4246

$ws.Range("H89").Value = 1866.6666
$ws.Range("J89").Value = 2000
$ws.Range("L89").Value = 10000
$ws.Range("N89").Value = -21232

$ws.Range("H113").Value = 6918
$ws.Range("I113").Value = 8002
$ws.Range("K113").Value = 8002
$ws.Range("M113").Value = -4748

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3208.1
$ws.Range("J45").Value = 4020.8333
$ws.Range("L45").Value = 4020.8333
$ws.Range("N45").Value = -4774.8333

$ws.Range("H102").Value = 7625
$ws.Range("I102").Value = 500
$ws.Range("K102").Value = 500
$ws.Range("M102").Value = 1122

$ws.Range("H132").Value = 2186.875
$ws.Range("I132").Value = 2186.875
$ws.Range("K132").Value = 6560.625
$ws.Range("M132").Value = -4030.625

$ws.Range("H137").Value = 74999.336
$ws.Range("J137").Value = 74999.336
$ws.Range("L137").Value = 74999.336
$ws.Range("N137").Value = -85199.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 749.44446
$ws.Range("I94").Value = 749.44446
$ws.Range("K94").Value = 749.44446
$ws.Range("M94").Value = -298.44446

$ws.Range("H103").Value = 26862.5
$ws.Range("J103").Value = 25816.666
$ws.Range("L103").Value = 25816.666
$ws.Range("N103").Value = -28160.666

$ws.Range("H105").Value = 1328.8667
$ws.Range("I105").Value = 1245.0834
$ws.Range("J105").Value = 1664
$ws.Range("K105").Value = 1245.0834
$ws.Range("L105").Value = 1664
$ws.Range("M105").Value = 501.9166
$ws.Range("N105").Value = -5158

$ws.Range("H134").Value = 1089.25
$ws.Range("I134").Value = 1089.25
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3267.75
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()

$ws.Range("H140").Value = 88260
$ws.Range("J140").Value = 88260
$ws.Range("L140").Value = 88260
$ws.Range("N140").Value = -98620

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2439.5557
$ws.Range("I58").Value = 1155.2307
$ws.Range("J58").Value = 5778.8
$ws.Range("K58").Value = 1155.2307
$ws.Range("L58").Value = 5778.8
$ws.Range("M58").Value = -952.2307000000001
$ws.Range("N58").Value = -6184.8

$ws.Range("H122").Value = 1065
$ws.Range("I122").Value = 1078
$ws.Range("K122").Value = 3234
$ws.Range("M122").Value = -784

$ws.Range("H132").Value = 1951.5385
$ws.Range("I132").Value = 1951.5385
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5854.6155
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 2916.25
$ws.Range("I134").Value = 1975.7142
$ws.Range("J134").Value = 9500
$ws.Range("K134").Value = 5927.142599999999
$ws.Range("L134").Value = 28500
$ws.Range("M134").Value = -3392.142599999999
$ws.Range("N134").Value = -33570

$ws.Range("H136").Value = 2439.5557
$ws.Range("I136").Value = 1155.2307
$ws.Range("J136").Value = 5778.8
$ws.Range("K136").Value = 3465.6921
$ws.Range("L136").Value = 17336.4
$ws.Range("M136").Value = -915.6921000000002
$ws.Range("N136").Value = -22436.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 150
$ws.Range("I8").Value = 150
$ws.Range("K8").Value = 450
$ws.Range("M8").Value = -311

$ws.Range("H113").Value = 736.4
$ws.Range("I113").Value = 477.33334
$ws.Range("J113").Value = 1125
$ws.Range("K113").Value = 1432.00002
$ws.Range("L113").Value = 3375
$ws.Range("M113").Value = 737.9999800000001
$ws.Range("N113").Value = -7715

$ws.Range("H131").Value = 882.375
$ws.Range("I131").Value = 911.2857
$ws.Range("J131").Value = 680
$ws.Range("K131").Value = 2733.8571
$ws.Range("L131").Value = 2040
$ws.Range("M131").Value = 2306.1429
$ws.Range("N131").Value = -12120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1850.3478
$ws.Range("I102").Value = 1545.7646
$ws.Range("J102").Value = 2713.3333
$ws.Range("K102").Value = 1545.7646
$ws.Range("L102").Value = 2713.3333
$ws.Range("M102").Value = 76.23540000000003
$ws.Range("N102").Value = -5957.3333

$ws.Range("H126").Value = 2602.3845
$ws.Range("I126").Value = 1994.25
$ws.Range("J126").Value = 9900
$ws.Range("K126").Value = 5982.75
$ws.Range("L126").Value = 29700
$ws.Range("M126").Value = -3512.75
$ws.Range("N126").Value = -34640

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8192.25
$ws.Range("I7").Value = 7773.1665
$ws.Range("K7").Value = 7773.1665
$ws.Range("M7").Value = -7661.1665

$ws.Range("H22").Value = 2093.75
$ws.Range("J22").Value = 2250
$ws.Range("L22").Value = 2250
$ws.Range("N22").Value = -2840

$ws.Range("H27").Value = 2093.75
$ws.Range("J27").Value = 2250
$ws.Range("L27").Value = 2250
$ws.Range("N27").Value = -2464

$ws.Range("H40").Value = 6103.6665
$ws.Range("I40").Value = 4989.4287
$ws.Range("K40").Value = 4989.4287
$ws.Range("M40").Value = -4853.4287

$ws.Range("H61").Value = 6221
$ws.Range("I61").Value = 4886.75
$ws.Range("K61").Value = 4886.75
$ws.Range("M61").Value = -4684.75

$ws.Range("H113").Value = 6221
$ws.Range("I113").Value = 4886.75
$ws.Range("K113").Value = 4886.75
$ws.Range("M113").Value = -2716.75

$ws.Range("H122").Value = 2665.5
$ws.Range("I122").Value = 2498.6
$ws.Range("K122").Value = 7495.799999999999
$ws.Range("M122").Value = -5045.799999999999

$ws.Range("H126").Value = 8192.25
$ws.Range("I126").Value = 7773.1665
$ws.Range("K126").Value = 23319.4995
$ws.Range("M126").Value = -20849.4995

$ws.Range("H132").Value = 998.3333
$ws.Range("I132").Value = 998.3333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2994.9999
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 1312.5
$ws.Range("J136").Value = 2000
$ws.Range("L136").Value = 6000
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8073
$ws.Range("I62").Value = 4800.5
$ws.Range("K62").Value = 4800.5
$ws.Range("M62").Value = -4176.5

$ws.Range("H65").Value = 8073
$ws.Range("I65").Value = 4800.5
$ws.Range("K65").Value = 24002.5
$ws.Range("M65").Value = -20882.5

$ws.Range("H81").Value = 351
$ws.Range("I81").Value = 351
$ws.Range("K81").Value = 702
$ws.Range("M81").Value = 359

$ws.Range("H84").Value = 351
$ws.Range("I84").Value = 351
$ws.Range("K84").Value = 3510
$ws.Range("M84").Value = 1794

$ws.Range("H132").Value = 1150.1333
$ws.Range("I132").Value = 1150.1333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3450.3999
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 3065.3809
$ws.Range("I136").Value = 2587.889
$ws.Range("K136").Value = 7763.667
$ws.Range("M136").Value = -5213.667

$ws.Range("H141").Value = 183570
$ws.Range("J141").Value = 120284
$ws.Range("L141").Value = 120284
$ws.Range("N141").Value = -130644
